# Update the crawl timestamp for every data row, and mark a handful of
# products as "Online kein Bestand" (out of stock online) in their
# productAriaLabel (column M) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2023-01-27 06:49:20"
$newTimestamp = "2023-01-27 12:55:52"

# Data rows run from row 2 through row 393 (row 1 is the header).
$lastRow = 393

# Map of row -> new productAriaLabel text (column M)
$labelUpdates = @{
    24  = "Betty Bossi Frischback Buttergipfel IP-Suisse - Online kein Bestand 2.60 Schweizer Franken"
    67  = "Betty Bossi Frischback Butterweggli IP-Suisse 6x50g - Online kein Bestand 3.10 Schweizer Franken"
    90  = "Betty Bossi Frischback Semmeli IP-Suisse - Online kein Bestand 2.40 Schweizer Franken"
    198 = "Pasquier Pitch Schokolade 8 Stück - Online kein Bestand 4.50 Schweizer Franken"
    216 = "Old El Paso Tortilla glutenfrei 6 Stück - Online kein Bestand 20% ab 2 Aktion 5.95 Schweizer Franken"
    228 = "Betty Bossi Spitzbube 2x  80g - Online kein Bestand 4.40 Schweizer Franken"
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # column O = timestamp
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }

    if ($labelUpdates.ContainsKey($row)) {
        $ws.Cells.Item($row, 13).Value = $labelUpdates[$row]
    }
}
